$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Korrekturen an bestehenden Stundenwerten ---
$ws.Range("B12").Value = 4
$ws.Range("B23").Value = 7
$ws.Range("B25").Value = 9
$ws.Range("B28").Value = 10
$ws.Range("B31").Value = 11

# --- Neue Zeile 32: Erfassung von gestern (21.01.2025) ---
# Formatierung (Datum in A, Zahl in B) von Zeile 31 uebernehmen, damit
# dieselben Formatvorlagen (cellXfs) wiederverwendet werden.
$ws.Range("A31:D31").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A32").Value = 45678
$ws.Range("B32").Value = 13
$ws.Range("C32").Value = "Aufgaben"
$ws.Range("D32").Value = "Lastenheft und GUI"

# --- Summenformel auf die neue Zeile erweitern ---
# (Zelle vorher leeren, damit die alte Formel/Abhaengigkeit sauber
# durch die neue ersetzt wird und der Cache korrekt neu aufgebaut wird.)
$ws.Range("F6").Value = ""
$ws.Range("F6").Formula = "=SUM(B7:B32)"

# --- Ansicht: Zoom anpassen und neue Zeile selektieren ---
$ws.Activate()
$null = $ws.Range("F21").Select()
$excel.ActiveWindow.Zoom = 130
